$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the analysis timestamp string in A2
$ws.Range("A2").Value = "2025-05-22 15:29:05"

# Update the metric values in row 2
$ws.Range("B2").Value = 16217
$ws.Range("C2").Value = 11724
$ws.Range("D2").Value = 72.29450576555466
$ws.Range("E2").Value = 2247
$ws.Range("F2").Value = 13.85583030153543
$ws.Range("G2").Value = 3106
$ws.Range("H2").Value = 19.15274095085404
$ws.Range("I2").Value = 9561
$ws.Range("J2").Value = 58.95665042856262
$ws.Range("K2").Value = 3080559.16
$ws.Range("L2").Value = 3550
$ws.Range("M2").Value = 21.89060862058334
$ws.Range("N2").Value = 1136656.43
$ws.Range("O2").Value = 4672
$ws.Range("P2").Value = 28.80927421841278
$ws.Range("Q2").Value = 489961.1
$ws.Range("R2").Value = 3544
$ws.Range("S2").Value = 21.85361040883024
$ws.Range("T2").Value = 3449
$ws.Range("U2").Value = 21.26780538940618
$ws.Range("V2").Value = 2436238.06
$ws.Range("W2").Value = 2306
$ws.Range("X2").Value = 14.2196460504409
$ws.Range("Y2").Value = 1440
$ws.Range("Z2").Value = 8.879570820743664
$ws.Range("AA2").Value = 154360
$ws.Range("AB2").Value = 806
$ws.Range("AC2").Value = 4.970093112166245
$ws.Range("AD2").Value = 16229
$ws.Range("AE2").Value = 10598
$ws.Range("AF2").Value = 65.30285291761662
$ws.Range("AG2").Value = 5631
$ws.Range("AH2").Value = 34.69714708238338
$ws.Range("AI2").Value = 466
$ws.Range("AJ2").Value = 864
$ws.Range("AK2").Value = 1366
$ws.Range("AL2").Value = 17.28486646884273
$ws.Range("AM2").Value = 32.04747774480713
$ws.Range("AN2").Value = 50.66765578635015
$ws.Range("AO2").Value = 1381340.64
$ws.Range("AP2").Value = 259487.28
$ws.Range("AQ2").Value = 86475.97
$ws.Range("AR2").Value = 79.97091004061826
$ws.Range("AS2").Value = 15.02267675666498
$ws.Range("AT2").Value = 5.006413202716749
$ws.Range("AU2").Value = 46.24280782508631
$ws.Range("AV2").Value = 233.1918505942275
$ws.Range("AW2").Value = 558.110151187905

$wb.Save()
